$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1. CORE COMPETENCIES: collapse the three long detail paragraphs into a
#    single short summary line (category names only, bullet-separated).
# -----------------------------------------------------------------------
# Paragraph 6 = "Product Management & Strategy: ..."      -> rewrite text
# Paragraph 7 = "Technical Product Development: ..."      -> remove
# Paragraph 8 = "Platform & Infrastructure: ..."           -> remove
$bullet = [char]0x2022
$summaryPara = $d.Paragraphs.Item(6)
$summaryPara.Range.Text = "Product Management & Strategy $bullet Technical Product Development $bullet Platform & Infrastructure"

# Paragraph indices shift down by one after each delete, so paragraph 7
# is removed twice (it is the old paragraph 7, then the old paragraph 8).
$d.Paragraphs.Item(7).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()

# -----------------------------------------------------------------------
# 2. Insert a new "TECHNICAL SKILLS" section (one Heading2 + three detail
#    paragraphs) immediately before the closing "For a more detailed..."
#    paragraph at the end of the document.
# -----------------------------------------------------------------------
$closingText = "For a more detailed, full description of my experience, please visit my LinkedIn and Personal Site."

$closingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pText = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($pText -eq $closingText) {
        $closingIndex = $i
        break
    }
}

$closingPara = $d.Paragraphs.Item($closingIndex)
$closingPara.Range.InsertParagraphBefore()
$closingPara.Range.InsertParagraphBefore()
$closingPara.Range.InsertParagraphBefore()
$closingPara.Range.InsertParagraphBefore()

$headingPara = $d.Paragraphs.Item($closingIndex)
$headingPara.Style = "Heading2"
$headingPara.Range.Text = "TECHNICAL SKILLS"

$skill1 = $d.Paragraphs.Item($closingIndex + 1)
$skill1.Range.Text = "PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development"

$skill2 = $d.Paragraphs.Item($closingIndex + 2)
$skill2.Range.Text = "TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; API Development"

$skill3 = $d.Paragraphs.Item($closingIndex + 3)
$skill3.Range.Text = "PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Security & Compliance"

Write-Output "TECHNICAL SKILLS section inserted at paragraph $closingIndex; CORE COMPETENCIES collapsed."
